$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 6
$ws_ALC.Range("H6").Value = 590
$ws_ALC.Range("I6").Value = 590
$ws_ALC.Range("J6").Value = 0
$ws_ALC.Range("K6").Value = 1770
$ws_ALC.Range("L6").Value = 0
$ws_ALC.Range("M6").Value = -1658
$ws_ALC.Range("N6").ClearContents()

# ALC row 43
$ws_ALC.Range("H43").Value = 0
$ws_ALC.Range("I43").Value = 0
$ws_ALC.Range("J43").Value = 0
$ws_ALC.Range("K43").Value = 0
$ws_ALC.Range("L43").Value = 0
$ws_ALC.Range("M43").ClearContents()
$ws_ALC.Range("N43").ClearContents()

# ALC row 62
$ws_ALC.Range("H62").Value = 2285.7144
$ws_ALC.Range("I62").Value = 1800
$ws_ALC.Range("K62").Value = 1800
$ws_ALC.Range("M62").Value = -1176

# ALC row 65
$ws_ALC.Range("H65").Value = 2285.7144
$ws_ALC.Range("I65").Value = 1800
$ws_ALC.Range("K65").Value = 9000
$ws_ALC.Range("M65").Value = -5880

# ALC row 111
$ws_ALC.Range("H111").Value = 999.2857
$ws_ALC.Range("I111").Value = 999.2857
$ws_ALC.Range("K111").Value = 2997.8571
$ws_ALC.Range("M111").Value = 69.14289999999983

# ALC row 129
$ws_ALC.Range("H129").Value = 1005.44446
$ws_ALC.Range("I129").Value = 543.55554
$ws_ALC.Range("J129").Value = 1097.8223
$ws_ALC.Range("K129").Value = 1630.66662
$ws_ALC.Range("L129").Value = 3293.4669
$ws_ALC.Range("M129").Value = 3369.33338
$ws_ALC.Range("N129").Value = -13293.4669

# ALC row 137
$ws_ALC.Range("H137").Value = 1579.7797
$ws_ALC.Range("I137").Value = 1364.2858
$ws_ALC.Range("K137").Value = 4092.8574
$ws_ALC.Range("M137").Value = -1542.8574

# ARM row 34
$ws_ARM.Range("H34").Value = 0
$ws_ARM.Range("J34").Value = 0
$ws_ARM.Range("L34").Value = 0
$ws_ARM.Range("N34").ClearContents()

# ARM row 110
$ws_ARM.Range("H110").Value = 1719.1666
$ws_ARM.Range("I110").Value = 1577.6666
$ws_ARM.Range("J110").Value = 2002.1666
$ws_ARM.Range("K110").Value = 1577.6666
$ws_ARM.Range("L110").Value = 2002.1666
$ws_ARM.Range("M110").Value = 467.3334
$ws_ARM.Range("N110").Value = -6092.1666

# ARM row 133
$ws_ARM.Range("H133").Value = 41828.418
$ws_ARM.Range("I133").Value = 20000
$ws_ARM.Range("J133").Value = 43812.816
$ws_ARM.Range("K133").Value = 20000
$ws_ARM.Range("L133").Value = 43812.816
$ws_ARM.Range("M133").Value = -17470
$ws_ARM.Range("N133").Value = -48872.816

# BSM row 11
$ws_BSM.Range("H11").Value = 102.25
$ws_BSM.Range("I11").Value = 102.25
$ws_BSM.Range("J11").Value = 0
$ws_BSM.Range("K11").Value = 102.25
$ws_BSM.Range("L11").Value = 0
$ws_BSM.Range("M11").Value = 37.75
$ws_BSM.Range("N11").ClearContents()

# BSM row 20
$ws_BSM.Range("H20").Value = 30371.97
$ws_BSM.Range("I20").Value = 43351.957
$ws_BSM.Range("J20").Value = 2052
$ws_BSM.Range("K20").Value = 43351.957
$ws_BSM.Range("L20").Value = 2052
$ws_BSM.Range("M20").Value = -43104.957
$ws_BSM.Range("N20").Value = -2546

# BSM row 106
$ws_BSM.Range("H106").Value = 53000
$ws_BSM.Range("J106").Value = 53000
$ws_BSM.Range("L106").Value = 53000
$ws_BSM.Range("N106").Value = -55524

# CRP row 31
$ws_CRP.Range("H31").Value = 3596.4783
$ws_CRP.Range("I31").Value = 2998.182
$ws_CRP.Range("J31").Value = 4144.9165
$ws_CRP.Range("K31").Value = 2998.182
$ws_CRP.Range("L31").Value = 4144.9165
$ws_CRP.Range("M31").Value = -2703.182
$ws_CRP.Range("N31").Value = -4734.9165

# CRP row 34
$ws_CRP.Range("H34").Value = 3596.4783
$ws_CRP.Range("I34").Value = 2998.182
$ws_CRP.Range("J34").Value = 4144.9165
$ws_CRP.Range("K34").Value = 2998.182
$ws_CRP.Range("L34").Value = 4144.9165
$ws_CRP.Range("M34").Value = -2796.182
$ws_CRP.Range("N34").Value = -4548.9165

# CUL row 5
$ws_CUL.Range("H5").Value = 1440.9
$ws_CUL.Range("I5").Value = 1799.8
$ws_CUL.Range("J5").Value = 1082
$ws_CUL.Range("K5").Value = 5399.4
$ws_CUL.Range("L5").Value = 3246
$ws_CUL.Range("M5").Value = -5287.4
$ws_CUL.Range("N5").Value = -3470

# CUL row 11
$ws_CUL.Range("H11").Value = 229.13333
$ws_CUL.Range("I11").Value = 118.375
$ws_CUL.Range("J11").Value = 355.7143
$ws_CUL.Range("K11").Value = 355.125
$ws_CUL.Range("L11").Value = 1067.1429
$ws_CUL.Range("M11").Value = -215.125
$ws_CUL.Range("N11").Value = -1347.1429

# CUL row 42
$ws_CUL.Range("H42").Value = 1496.5
$ws_CUL.Range("I42").Value = 993
$ws_CUL.Range("K42").Value = 2979
$ws_CUL.Range("M42").Value = -2445

# CUL row 131
$ws_CUL.Range("H131").Value = 883.75
$ws_CUL.Range("J131").Value = 891.59186
$ws_CUL.Range("L131").Value = 2674.77558
$ws_CUL.Range("N131").Value = -12754.77558

# CUL row 132
$ws_CUL.Range("H132").Value = 2216.44
$ws_CUL.Range("I132").Value = 1200.091
$ws_CUL.Range("J132").Value = 3015
$ws_CUL.Range("K132").Value = 10800.819
$ws_CUL.Range("L132").Value = 27135
$ws_CUL.Range("M132").Value = -8270.819
$ws_CUL.Range("N132").Value = -32195

# CUL row 135
$ws_CUL.Range("H135").Value = 1440.9
$ws_CUL.Range("I135").Value = 1799.8
$ws_CUL.Range("J135").Value = 1082
$ws_CUL.Range("K135").Value = 16198.2
$ws_CUL.Range("L135").Value = 9738
$ws_CUL.Range("M135").Value = -13663.2
$ws_CUL.Range("N135").Value = -14808

# CUL row 137
$ws_CUL.Range("H137").Value = 37039760
$ws_CUL.Range("I137").Value = 1107.5
$ws_CUL.Range("J137").Value = 66670680
$ws_CUL.Range("K137").Value = 3322.5
$ws_CUL.Range("L137").Value = 200012040
$ws_CUL.Range("M137").Value = 1777.5
$ws_CUL.Range("N137").Value = -200022240

# GSM row 13
$ws_GSM.Range("H13").Value = 296.66666
$ws_GSM.Range("I13").Value = 296.66666
$ws_GSM.Range("K13").Value = 296.66666
$ws_GSM.Range("M13").Value = -157.66666

# GSM row 92
$ws_GSM.Range("H92").Value = 4076.375
$ws_GSM.Range("J92").Value = 4076.375
$ws_GSM.Range("L92").Value = 4076.375
$ws_GSM.Range("N92").Value = -7820.375

# GSM row 114
$ws_GSM.Range("H114").Value = 19900
$ws_GSM.Range("J114").Value = 19900
$ws_GSM.Range("L114").Value = 19900
$ws_GSM.Range("N114").Value = -28578

# LTW row 100
$ws_LTW.Range("H100").Value = 5873.769
$ws_LTW.Range("I100").Value = 8020
$ws_LTW.Range("J100").Value = 2439.8
$ws_LTW.Range("K100").Value = 8020
$ws_LTW.Range("L100").Value = 2439.8
$ws_LTW.Range("M100").Value = -7479
$ws_LTW.Range("N100").Value = -3521.8

# WVR row 21
$ws_WVR.Range("H21").Value = 19613.6
$ws_WVR.Range("J21").Value = 12017
$ws_WVR.Range("L21").Value = 12017
$ws_WVR.Range("N21").Value = -12487

# WVR row 35
$ws_WVR.Range("H35").Value = 19613.6
$ws_WVR.Range("J35").Value = 12017
$ws_WVR.Range("L35").Value = 12017
$ws_WVR.Range("N35").Value = -12597

# WVR row 63
$ws_WVR.Range("H63").Value = 108459.6
$ws_WVR.Range("J63").Value = 108459.6
$ws_WVR.Range("L63").Value = 108459.6
$ws_WVR.Range("N63").Value = -109707.6

# WVR row 66
$ws_WVR.Range("H66").Value = 108459.6
$ws_WVR.Range("J66").Value = 108459.6
$ws_WVR.Range("L66").Value = 325378.8
$ws_WVR.Range("N66").Value = -331618.8

# WVR row 82
$ws_WVR.Range("H82").Value = 0
$ws_WVR.Range("J82").Value = 0
$ws_WVR.Range("L82").Value = 0
$ws_WVR.Range("N82").ClearContents()

# WVR row 85
$ws_WVR.Range("H85").Value = 0
$ws_WVR.Range("J85").Value = 0
$ws_WVR.Range("L85").Value = 0
$ws_WVR.Range("N85").ClearContents()

# WVR row 107
$ws_WVR.Range("H107").Value = 528.8182
$ws_WVR.Range("I107").Value = 411.07693
$ws_WVR.Range("K107").Value = 1233.23079
$ws_WVR.Range("M107").Value = 686.7692099999999

# WVR row 122
$ws_WVR.Range("H122").Value = 14708411
$ws_WVR.Range("I122").Value = 17858700
$ws_WVR.Range("J122").Value = 7063.3335
$ws_WVR.Range("K122").Value = 53576100
$ws_WVR.Range("L122").Value = 21190.0005
$ws_WVR.Range("M122").Value = -53573650
$ws_WVR.Range("N122").Value = -26090.0005

# WVR row 123
$ws_WVR.Range("H123").Value = 42071
$ws_WVR.Range("J123").Value = 42071
$ws_WVR.Range("L123").Value = 42071
$ws_WVR.Range("N123").Value = -51871

# WVR row 136
$ws_WVR.Range("H136").Value = 1328.6029
$ws_WVR.Range("I136").Value = 1303.88
$ws_WVR.Range("J136").Value = 1397.2778
$ws_WVR.Range("K136").Value = 3911.64
$ws_WVR.Range("L136").Value = 4191.8334
$ws_WVR.Range("M136").Value = -1361.64
$ws_WVR.Range("N136").Value = -9291.8334
